$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rewrite the existing data rows (2-6) in place so each cell keeps the
#    per-cell style that was already on it (A=Keyword col, B=Locator col,
#    C=Data col).
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,1).Value = "clickByXpath"
$ws.Cells.Item(2,2).Value = "/html/body/div[2]/div[1]/div/div[1]/div[2]/a[3]"
$ws.Cells.Item(2,3).Value = ""

$ws.Cells.Item(3,1).Value = "clickByXpath"
$ws.Cells.Item(3,2).Value = "/html/body/div[20]/div/div[2]/div[2]/div[2]/form/div[1]/input"
$ws.Cells.Item(3,3).Value = ""

$ws.Cells.Item(4,1).Value = "enterByXpath"
$ws.Cells.Item(4,2).Value = "/html/body/div[20]/div/div[2]/div[2]/div[2]/form/div[1]/input"
$ws.Cells.Item(4,3).Value = "rohithkumar90@ymail.com"

$ws.Cells.Item(5,1).Value = "clickByXpath"
$ws.Cells.Item(5,2).Value = "/html/body/div[20]/div/div[2]/div[2]/div[2]/form/div[2]/input"
$ws.Cells.Item(5,3).Value = ""

$ws.Cells.Item(6,1).Value = "enterByXpath"
$ws.Cells.Item(6,2).Value = "/html/body/div[20]/div/div[2]/div[2]/div[2]/form/div[2]/input"
$ws.Cells.Item(6,3).Value = "rohith270419909940096410abcd1990"

# ---------------------------------------------------------------------------
# 2) Append new rows 7-10. Clone row 6's formatting down across them first
#    (so the new rows pick up the same A/B/C column styles used throughout
#    the table), then fill in the new keyword/locator values. Row 10 is
#    left blank - it only carries the inherited formatting.
# ---------------------------------------------------------------------------
$ws.Range("A6:C6").Copy()
$ws.Range("A7:C10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Cells.Item(7,1).Value = "clickByXpath"
$ws.Cells.Item(7,2).Value = "/html/body/div[20]/div/div[2]/div[2]/div[2]/form/input[1]"

$ws.Cells.Item(8,1).Value = "clickByXpath"
$ws.Cells.Item(8,2).Value = "/html/body/div[2]/div[1]/div[2]/div/div/div[3]/ul/li[1]/a"

$ws.Cells.Item(9,1).Value = "clickByXpath"
$ws.Cells.Item(9,2).Value = "/html/body/div[2]/div[1]/div[2]/div/div/div[3]/ul/li[1]/div/ul/li[13]/button"

# Row 10 stays empty (A10/B10/C10 keep the copied format only).

# ---------------------------------------------------------------------------
# 3) Turn the e-mail address in C4 into a real hyperlink (adds the
#    built-in "Hyperlink" cell style / underlined theme font automatically).
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("C4"), "mailto:rohithkumar90@ymail.com")

# ---------------------------------------------------------------------------
# 4) Column B needs to widen (no longer auto "best fit") to comfortably show
#    the long xpath locator strings.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 76.6

# ---------------------------------------------------------------------------
# 5) Update the view: the active selection moves to C8 and the window is
#    scrolled down one row so row 2 sits at the top of the viewport.
# ---------------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 2
$ws.Range("C8").Select()
